# Applies the crypto price/volume refresh described in the commit diff.
# D (Price) and E (Volume 1h) columns hold numeric-looking text stored as
# plain strings in the sheet, so each new value is written with a leading
# apostrophe to force Excel to keep it as text (preventing float rounding
# such as 325.86 -> 325.86000000000001, and keeping the literal '%' text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'325.86"
$ws.Range('E2').Value = "'-2.73%"
$ws.Range('D3').Value = "'44.22"
$ws.Range('E3').Value = "'0.86%"
$ws.Range('D4').Value = "'5.587"
$ws.Range('E4').Value = "'-2.88%"
$ws.Range('D5').Value = "'0.08046"
$ws.Range('E5').Value = "'-3.91%"
$ws.Range('D6').Value = "'4.298"
$ws.Range('E6').Value = "'-4.83%"
$ws.Range('D7').Value = "'1.886"
$ws.Range('E7').Value = "'-3.64%"
$ws.Range('D8').Value = "'2.647"
$ws.Range('E8').Value = "'-8.08%"
$ws.Range('D9').Value = "'0.9445"
$ws.Range('E9').Value = "'-0.45%"
$ws.Range('D10').Value = "'0.1160"
$ws.Range('E10').Value = "'-6.95%"
$ws.Range('D11').Value = "'0.1835"
$ws.Range('E11').Value = "'-6.87%"
$ws.Range('D12').Value = "'0.09808"
$ws.Range('E12').Value = "'-5.77%"
$ws.Range('D13').Value = "'0.04265"
$ws.Range('E13').Value = "'-6.42%"
$ws.Range('D14').Value = "'0.1065"
$ws.Range('E14').Value = "'-0.16%"
$ws.Range('D15').Value = "'0.001271"
$ws.Range('E15').Value = "'-2.52%"
$ws.Range('D16').Value = "'0.04216"
$ws.Range('E16').Value = "'-4.77%"
$ws.Range('D17').Value = "'0.005952"
$ws.Range('E17').Value = "'0.73%"
$ws.Range('E18').Value = "'3.17%"
$ws.Range('D19').Value = "'0.3496"
$ws.Range('E19').Value = "'-0.27%"
$ws.Range('D20').Value = "'8.333"
$ws.Range('E20').Value = "'-4.39%"
$ws.Range('D21').Value = "'0.1379"
$ws.Range('E21').Value = "'1.21%"
$ws.Range('E22').Value = "'0.70%"
$ws.Range('E23').Value = "'-0.85%"
$ws.Range('D24').Value = "'0.004508"
$ws.Range('E24').Value = "'3.63%"
$ws.Range('D25').Value = "'0.0001262"
$ws.Range('E25').Value = "'-0.05%"
$ws.Range('D26').Value = "'0.0003993"
$ws.Range('E26').Value = "'0.00%"
$ws.Range('E38').Value = "'-6.66%"
$ws.Range('D39').Value = "'0.05438"
$ws.Range('E39').Value = "'-10.48%"
$ws.Range('D40').Value = "'0.007587"
$ws.Range('E40').Value = "'-4.17%"
$ws.Range('D41').Value = "'0.1394"
$ws.Range('E41').Value = "'-2.25%"
$ws.Range('D42').Value = "'0.007330"
$ws.Range('E42').Value = "'-18.33%"
$ws.Range('D43').Value = "'0.002019"
$ws.Range('E43').Value = "'-5.91%"
$ws.Range('D44').Value = "'0.008812"
$ws.Range('E44').Value = "'-13.03%"
$ws.Range('D45').Value = "'0.00006910"
$ws.Range('E45').Value = "'-4.95%"
$ws.Range('D46').Value = "'0.00000000751"
$ws.Range('E46').Value = "'0.00%"
$ws.Range('B47').Value = 'CoinbaseStockToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range('D47').Value = "'0.002272"
$ws.Range('E47').Value = "'0.00%"
$ws.Range('B48').Value = 'BOLO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range('D48').Value = "'0.003781"
$ws.Range('E48').Value = "'18.56%"
$ws.Range('D49').Value = "'0.00002102"
$ws.Range('E49').Value = "'0.00%"
$ws.Range('D50').Value = "'0.0002002"
$ws.Range('E50').Value = "'0.00%"
